$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose displayed text changes. Values are kept as plain text
# (matching the original inline-string cell type) rather than letting
# Excel auto-convert numeric-looking / percentage-looking strings into
# numbers.
$updates = [ordered]@{
    "D2" = "329.30"
    "E2" = "1.48%"
    "D3" = "41.34"
    "E3" = "4.38%"
    "D4" = "5.620"
    "E4" = "-0.17%"
    "D5" = "0.08217"
    "E5" = "2.45%"
    "E6" = "1.53%"
    "D7" = "2.002"
    "E7" = "-0.77%"
    "D8" = "4.501"
    "E8" = "0.22%"
    "E9" = "1.93%"
    "D10" = "0.9226"
    "E10" = "0.02%"
    "E11" = "2.85%"
    "D12" = "0.1960"
    "E12" = "-0.03%"
    "E13" = "1.65%"
    "D14" = "0.03847"
    "E14" = "7.71%"
    "D15" = "0.1059"
    "E15" = "1.03%"
    "E16" = "0.52%"
    "D17" = "0.006229"
    "E17" = "-0.23%"
    "D19" = "3.449"
    "E19" = "2.94%"
    "E20" = "-0.01%"
    "D21" = "8.269"
    "E21" = "-5.15%"
    "D22" = "0.1366"
    "E22" = "0.89%"
    "D23" = "0.2662"
    "E23" = "6.17%"
    "D24" = "0.04408"
    "E24" = "0.67%"
    "D25" = "0.001258"
    "E25" = "-0.42%"
    "E26" = "-6.31%"
    "E27" = "-2.49%"
    "D39" = "0.02747"
    "E39" = "10.40%"
    "D40" = "0.05449"
    "E40" = "2.44%"
    "D41" = "0.007965"
    "E41" = "7.21%"
    "D42" = "0.1421"
    "E42" = "1.27%"
    "D43" = "0.008940"
    "D44" = "0.002172"
    "E44" = "2.52%"
    "D45" = "0.01144"
    "E45" = "2.50%"
    "D46" = "0.00006776"
    "E46" = "1.23%"
    "D47" = "0.00000000751"
    "E47" = "-0.05%"
    "D48" = "0.003192"
    "E48" = "7.19%"
    "D49" = "0.002280"
    "E49" = "-0.08%"
    "D50" = "0.00002102"
    "E50" = "-0.05%"
    "D51" = "0.0002002"
    "E51" = "-0.05%"
}

foreach ($name in $updates.Keys) {
    $cell = $ws.Range($name)
    # Force a text format before assigning so Excel does not reinterpret
    # strings such as "329.30" or "1.48%" as numbers/percentages.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$name]
    # Restore the default style so no stray number-format index is left
    # attached to the cell (original cells carry no explicit style here).
    $cell.Style = "Normal"
}
